# Reg TCs for 811, 2805, 2648 stories
#
# Inserts a new "pop5" test-data block (Duplicate_PublicationID scenario)
# above the existing pop5 block at row 141, pushing the old pop5 block
# (which becomes "pop6") and the old pop6 block (which becomes "pop7")
# down by 5 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room for the new block: insert 5 blank rows at row 141.
#    Everything that used to start at row 141 (old pop5 block, old gap,
#    old pop6 row) shifts down to row 146 onward.
$ws.Rows("141:145").Insert()

# 2) Fill in the new "pop5" block (rows 141-144, row 145 stays blank as
#    the separator row before the next block, matching the sheet's
#    existing convention).
$ws.Range("A141").Value = "pop5"
$ws.Range("C141").Value = "LIVEHTA Automation - Test_NonOncology_Automation_1"
$ws.Range("E141").Value = "\Testdata\Non_Oncology\Templates\ImportPublications\Extraction sheet - Duplicate_PublicationID.xlsx"
$ws.Range("D141").Value = "Extraction sheet - Duplicate_PublicationID.xlsx"
$ws.Range("F141").Value = 13
$ws.Range("G141").Value = "Primary Publication " + [char]0x201C + "Feagan_NEJM_2013" + [char]0x201D + " is also present in Related Publications"

$ws.Range("A142").Value = "pop5"
$ws.Range("F142").Value = 14
$ws.Range("G142").Value = "Primary Publication " + [char]0x201C + "Feagan_NEJM_2013" + [char]0x201D + " is also present in Related Publications"

$ws.Range("A143").Value = "pop5"
$ws.Range("F143").Value = 22
$ws.Range("G143").Value = "Primary Publication " + [char]0x201C + "Motoya_PLOS_2019" + [char]0x201D + " is also present in Related Publications"

$ws.Range("A144").Value = "pop5"
$ws.Range("F144").Value = 23
$ws.Range("G144").Value = "Primary Publication " + [char]0x201C + "Motoya_PLOS_2019" + [char]0x201D + " is also present in Related Publications"

# 3) The block that used to be "pop5" (rows 141-151) is now at rows
#    146-156; relabel its Name column from pop5 to pop6.
for ($r = 146; $r -le 156; $r++) {
    $ws.Cells.Item($r, 1).Value = "pop6"
}

# 4) The single row that used to be "pop6" (row 153) is now at row 158;
#    relabel its Name column from pop6 to pop7.
$ws.Range("A158").Value = "pop7"

# 5) Column D best-fit width shrinks now that the widest new entry is
#    gone; get as close as the host allows to the recorded width.
$ws.Columns("D:D").ColumnWidth = 38.6

# 6) Leave the selection where the user ended up after the edit.
$ws.Range("A159").Select()
